# This script applies scheduled market-price / profit recalculation updates
# to the leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, columns H-N (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ)
# are refreshed with newly computed values; some cells are cleared entirely
# where no current data applies.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "ALC"; Sets = @{ "H2"=237.5; "I2"=237.5; "J2"=0; "K2"=237.5; "L2"=0; "M2"=-124.5 }; Clears = @("N2") }
    @{ Sheet = "ALC"; Sets = @{ "H32"=938.3333; "I32"=716; "K32"=716; "M32"=-390 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H33"=105.46667; "I33"=95.7; "K33"=95.7; "M33"=133.3 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H74"=3743.5557; "I74"=3379.2; "K74"=3379.2; "M74"=-2443.2 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H77"=3743.5557; "I77"=3379.2; "K77"=16896; "M77"=-12216 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H100"=1000; "I100"=1000; "K100"=1000; "M100"=-459 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H137"=1694.6086; "I137"=1255.4445; "J137"=1976.9286; "K137"=3766.3335; "L137"=5930.7858; "M137"=-1216.3335; "N137"=-11030.7858 }; Clears = @() }
    @{ Sheet = "ALC"; Sets = @{ "H138"=2859.36; "I138"=3029.25; "J138"=2179.8; "K138"=9087.75; "L138"=6539.400000000001; "M138"=-3947.75; "N138"=-16819.4 }; Clears = @() }
    @{ Sheet = "ARM"; Sets = @{ "H32"=4310.766; "I32"=2789.0278; "J32"=9291; "K32"=2789.0278; "L32"=9291; "M32"=-2502.0278; "N32"=-9865 }; Clears = @() }
    @{ Sheet = "ARM"; Sets = @{ "H45"=1826.6875; "I45"=875; "J45"=3050.2856; "K45"=875; "L45"=3050.2856; "M45"=-498; "N45"=-3804.2856 }; Clears = @() }
    @{ Sheet = "ARM"; Sets = @{ "H74"=499.66666; "I74"=499.66666; "K74"=499.66666; "M74"=374.33334 }; Clears = @() }
    @{ Sheet = "ARM"; Sets = @{ "H77"=499.66666; "I77"=499.66666; "K77"=2498.3333; "M77"=1869.6667 }; Clears = @() }
    @{ Sheet = "ARM"; Sets = @{ "H109"=75914.25; "J109"=75914.25; "L109"=75914.25; "N109"=-78688.25 }; Clears = @() }
    @{ Sheet = "BSM"; Sets = @{ "H29"=0; "I29"=0; "K29"=0 }; Clears = @("M29") }
    @{ Sheet = "BSM"; Sets = @{ "H39"=0; "J39"=0; "L39"=0 }; Clears = @("N39") }
    @{ Sheet = "CRP"; Sets = @{ "H16"=507.66666; "I16"=11; "J16"=756; "K16"=11; "L16"=756; "M16"=276; "N16"=-1330 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H41"=29000; "J41"=29000; "L41"=29000; "N41"=-29856 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H50"=18000; "J50"=18000; "L50"=18000; "N50"=-19250 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H51"=30283.334; "I51"=0; "J51"=30283.334; "K51"=0; "L51"=30283.334; "N51"=-31755.334 }; Clears = @("M51") }
    @{ Sheet = "CRP"; Sets = @{ "H59"=33266.668; "J59"=33266.668; "L59"=33266.668; "N59"=-35556.668 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H60"=11098.261; "J60"=11098.261; "L60"=11098.261; "N60"=-12120.261 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H61"=30283.334; "I61"=0; "J61"=30283.334; "K61"=0; "L61"=30283.334; "N61"=-30979.334 }; Clears = @("M61") }
    @{ Sheet = "CRP"; Sets = @{ "H94"=883.8333; "I94"=713; "J94"=1005.8571; "K94"=713; "L94"=1005.8571; "M94"=-262; "N94"=-1907.8571 }; Clears = @() }
    @{ Sheet = "CRP"; Sets = @{ "H113"=507.66666; "I113"=11; "J113"=756; "K113"=11; "L113"=756; "M113"=2159; "N113"=-5096 }; Clears = @() }
    @{ Sheet = "CUL"; Sets = @{ "H118"=2524.5; "I118"=715.6667; "J118"=4333.3335; "K118"=2147.0001; "L118"=13000.0005; "M118"=-904.0001000000002; "N118"=-15486.0005 }; Clears = @() }
    @{ Sheet = "CUL"; Sets = @{ "H130"=1006.6667; "I130"=1010; "J130"=1000; "K130"=3030; "L130"=3000; "M130"=1990; "N130"=-13040 }; Clears = @() }
    @{ Sheet = "GSM"; Sets = @{ "H20"=1801715.8; "J20"=12000; "L20"=12000; "N20"=-12490 }; Clears = @() }
    @{ Sheet = "GSM"; Sets = @{ "H24"=2509527.8; "I24"=10000000; "J24"=12703.667; "K24"=10000000; "L24"=12703.667; "M24"=-9999827; "N24"=-13049.667 }; Clears = @() }
    @{ Sheet = "GSM"; Sets = @{ "H102"=2154.862; "I102"=2127.611; "J102"=2199.4546; "K102"=2127.611; "L102"=2199.4546; "M102"=-505.6109999999999; "N102"=-5443.4546 }; Clears = @() }
    @{ Sheet = "GSM"; Sets = @{ "H126"=50295.477; "I126"=2867.5833; "J126"=113532.664; "K126"=8602.749899999999; "L126"=340597.992; "M126"=-6132.749899999999; "N126"=-345537.992 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H7"=9317; "J7"=9179.6; "L7"=9179.6; "N7"=-9403.6 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H46"=1385.7142; "I46"=1200; "K46"=1200; "M46"=-1012 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H100"=200; "I100"=200; "K100"=200; "M100"=341 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H126"=9317; "J126"=9179.6; "L126"=27538.8; "N126"=-32478.8 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H132"=2128.3157; "I132"=1150; "J132"=2243.4119; "K132"=3450; "L132"=6730.2357; "M132"=-920; "N132"=-11790.2357 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H136"=3569.7368; "I136"=2173.75; "J136"=4585; "K136"=6521.25; "L136"=13755; "M136"=-3971.25; "N136"=-18855 }; Clears = @() }
    @{ Sheet = "LTW"; Sets = @{ "H139"=48278.75; "J139"=48278.75; "L139"=48278.75; "N139"=-58558.75 }; Clears = @() }
    @{ Sheet = "WVR"; Sets = @{ "H21"=0; "J21"=0; "L21"=0 }; Clears = @("N21") }
    @{ Sheet = "WVR"; Sets = @{ "H24"=12000; "J24"=12000; "L24"=12000; "N24"=-12460 }; Clears = @() }
    @{ Sheet = "WVR"; Sets = @{ "H31"=0; "I31"=0; "J31"=0; "K31"=0; "L31"=0 }; Clears = @("M31", "N31") }
    @{ Sheet = "WVR"; Sets = @{ "H35"=0; "J35"=0; "L35"=0 }; Clears = @("N35") }
    @{ Sheet = "WVR"; Sets = @{ "H92"=28591.666; "J92"=28591.666; "L92"=28591.666; "N92"=-33583.666 }; Clears = @() }
    @{ Sheet = "WVR"; Sets = @{ "H132"=3960.2666; "I132"=3700.5833; "K132"=11101.7499; "M132"=-8571.749899999999 }; Clears = @() }
    @{ Sheet = "WVR"; Sets = @{ "H136"=3999.6667; "I136"=0; "J136"=3999.6667; "K136"=0; "L136"=11999.0001; "N136"=-17099.0001 }; Clears = @("M136") }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    foreach ($key in $u.Sets.Keys) {
        $ws.Range($key).Value = $u.Sets[$key]
    }
    foreach ($ref in $u.Clears) {
        $ws.Range($ref).ClearContents()
    }
}
